$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (A1:U1): _old -> _FV2404, _new -> _FV2410 -----------
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$headersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}

# Column K (11) holds "diff" and stays unchanged.

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# --- Turn the data range into an Excel Table (ListObject) ------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U61"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (split below row 1) ------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
